$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q: header "2020" in row 3 (same style pattern as P3=23)
$ws.Cells.Item(3, 17).Value = 2020
$ws.Cells.Item(3, 17).Style = $ws.Cells.Item(3, 16).Style

# Rows 4-13: fill with "-" to match the rest of the row, following the
# same style as column P for rows 5-13 (style 27/29); row 4 uses style 27
# (text) even though P4 uses a numeric style, because Q4 holds text "-".
$ws.Cells.Item(4, 17).Value = "-"
$ws.Cells.Item(4, 17).Style = $ws.Cells.Item(5, 16).Style

for ($r = 5; $r -le 12; $r++) {
    $ws.Cells.Item($r, 17).Value = "-"
    $ws.Cells.Item($r, 17).Style = $ws.Cells.Item($r, 16).Style
}

$ws.Cells.Item(13, 17).Value = "-"
$ws.Cells.Item(13, 17).Style = $ws.Cells.Item(13, 16).Style

# Update the active selection to match the authored worksheet view
$ws.Range("P17").Select()
